# Tuntikirjanpito.xlsx update
# "live versions for client and api are somewhat working, not great not terrible"
#
# Adds two new timesheet rows (68 and 69) to the "Taul1" sheet, describing
# the client login/signup work and the heroku live testing of the client+api,
# and extends the hours-total SUM formula to cover the new rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 68: 2 hours, client login/logout/signup work, tagged "client"
$ws.Range("B68").Value = 2
$ws.Range("C68").Value = "client login, logout ja signup pohja, cookie pohjainen permanent login lisätty "
$ws.Range("D68").Value = "client"

# Row 69: 1 hour, heroku live testing of client+api, tagged "api/client"
# (set D69 before C69 so the new shared strings are appended in the same
# order as in the target workbook: api/client, then the heroku description)
$ws.Range("B69").Value = 1
$ws.Range("D69").Value = "api/client"
$ws.Range("C69").Value = "heroku live testausta, client ja api livenä ja toimii jotenkin(login hidas), not great, not terrible"

# Extend the total-hours formula in row 75 to include the two new rows
$ws.Range("B75").Formula = "=SUM(B2:B69)"

# Update the selected cell to match where the editor last left off
$ws.Range("C69").Select()
